$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the two "Ticket" answers that were mis-marked ---
$ws.Range("J18").Value = "Ja"
$ws.Range("J19").Value = "Ja"

# --- Insert a new row so the list grows by one concert ---
# (old row 21 "Lizzy McAlpine" -> row 22, old row 22 "Only the Poets" -> row 23)
$ws.Rows("21:21").Insert()

# --- Row 20: replace the old "Chance Pena / Jonah Kagen" entry with "Chappell Roan" ---
$ws.Range("B20").ClearContents()
$ws.Range("A20").Value = "Chappell Roan"
$ws.Range("C20").Value = "Germany"
$ws.Range("D20").Value = "Berlin"
$ws.Range("E20").Value = "12000"
$ws.Range("F20").Value = "52.5308904332696"
$ws.Range("G20").Value = "13.451074492408171"
$ws.Range("H20").Value = "Velodrom"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("I20").Value = "10/23/2024"
$ws.Range("J20").Value = "Ja"
$ws.Range("K20").Value = "Daniela Wolfangel"
$ws.Range("K20").NumberFormat = "@"

# --- Row 21 (new row): "Peter Fox" concert at Görlitzer Park ---
$ws.Range("A21").Value = "Peter Fox"
$ws.Range("C21").Value = "Germany"
$ws.Range("D21").Value = "Berlin"
$ws.Range("E21").Value = "12000"
$ws.Range("F21").Value = "52.49678216757308"
$ws.Range("G21").Value = "13.437405995049385"
$ws.Range("H21").Value = "Görlitzer Park"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("I21").Value = "9/21/2024"
$ws.Range("J21").Value = "Nein"
$ws.Range("K21").Value = "Lisa Schmidt"
